$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the crypto symbol list refresh.
# Values are entered with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing inline-string cell format) instead of
# auto-converting numeric-looking / percentage-looking text into real numbers.
$ws.Range("D2").Value = "'306.92"
$ws.Range("E2").Value = "'0.25%"

$ws.Range("D3").Value = "'35.43"
$ws.Range("E3").Value = "'-2.46%"

$ws.Range("D4").Value = "'5.088"
$ws.Range("E4").Value = "'0.47%"

$ws.Range("D5").Value = "'0.08139"
$ws.Range("E5").Value = "'2.73%"

$ws.Range("D6").Value = "'1.986"
$ws.Range("E6").Value = "'-9.88%"

$ws.Range("D7").Value = "'7.912"
$ws.Range("E7").Value = "'-1.32%"

$ws.Range("D8").Value = "'2.876"
$ws.Range("E8").Value = "'9.27%"

$ws.Range("D9").Value = "'0.9246"
$ws.Range("E9").Value = "'-0.47%"

$ws.Range("D10").Value = "'0.1097"
$ws.Range("E10").Value = "'11.43%"

$ws.Range("D11").Value = "'0.1914"
$ws.Range("E11").Value = "'2.01%"

$ws.Range("D12").Value = "'0.09508"
$ws.Range("E12").Value = "'4.64%"

$ws.Range("D13").Value = "'0.03666"
$ws.Range("E13").Value = "'-1.20%"

$ws.Range("D14").Value = "'0.09910"
$ws.Range("E14").Value = "'-0.06%"

$ws.Range("D15").Value = "'0.001419"
$ws.Range("E15").Value = "'-1.02%"

$ws.Range("D16").Value = "'0.005814"
$ws.Range("E16").Value = "'3.19%"

$ws.Range("D17").Value = "'3.475"
$ws.Range("E17").Value = "'0.38%"

$ws.Range("D18").Value = "'4.138"

$ws.Range("E19").Value = "'1.44%"

$ws.Range("D20").Value = "'0.1313"
$ws.Range("E20").Value = "'-1.83%"

$ws.Range("D21").Value = "'5.105"
$ws.Range("E21").Value = "'-0.16%"

$ws.Range("D22").Value = "'0.2197"
$ws.Range("E22").Value = "'0.40%"

$ws.Range("D23").Value = "'0.04527"
$ws.Range("E23").Value = "'-0.53%"

$ws.Range("D24").Value = "'0.001228"
$ws.Range("E24").Value = "'-0.86%"

$ws.Range("D25").Value = "'0.004719"
$ws.Range("E25").Value = "'-1.33%"

$ws.Range("D26").Value = "'0.0001254"
$ws.Range("E26").Value = "'-3.49%"

$ws.Range("D27").Value = "'0.0004462"
$ws.Range("E27").Value = "'-5.83%"

$ws.Range("D39").Value = "'0.01939"
$ws.Range("E39").Value = "'0.85%"

$ws.Range("D40").Value = "'0.04839"
$ws.Range("E40").Value = "'-2.53%"

$ws.Range("D41").Value = "'0.007609"
$ws.Range("E41").Value = "'-2.57%"

$ws.Range("D42").Value = "'0.009673"
$ws.Range("E42").Value = "'23.92%"

$ws.Range("D43").Value = "'0.1369"
$ws.Range("E43").Value = "'-1.91%"

$ws.Range("D44").Value = "'0.002121"
$ws.Range("E44").Value = "'0.52%"

$ws.Range("D45").Value = "'0.01124"
$ws.Range("E45").Value = "'0.24%"

$ws.Range("D46").Value = "'0.00006519"
$ws.Range("E46").Value = "'4.71%"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.37%"

$ws.Range("E48").Value = "'24.03%"

$ws.Range("D49").Value = "'0.001304"
$ws.Range("E49").Value = "'-27.51%"

$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.37%"

$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.37%"

